# Re-ran the averaged-intensity computation including the new spiral
# orientation schemes. Gaussian-Quadrature moved up (re-computed with new
# values), three new "Spiral-*" rows were inserted after it, and the
# remaining rotation/hex-grid schemes were pushed down to make room.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.085420576528351
$ws.Range("D10").Value = 0.9436809254642285
$ws.Range("E10").Value = 0.986338913977008
$ws.Range("F10").Value = 0.975497084934187
$ws.Range("G10").Value = 1.085420576528351
$ws.Range("H10").Value = 0.9436809254642285
$ws.Range("I10").Value = 1.013507860391305
$ws.Range("J10").Value = 0.962453540013901
$ws.Range("K10").Value = 1.023791702434834
$ws.Range("L10").Value = 0.9516332046643771
$ws.Range("M10").Value = 1.085420576528351
$ws.Range("N10").Value = 0.9650099197206183
$ws.Range("O10").Value = 0.9977343752259437
$ws.Range("P10").Value = 0.992790476051024

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 1.060223973180004
$ws.Range("D11").Value = 0.8554390596845817
$ws.Range("E11").Value = 1.029087147232076
$ws.Range("F11").Value = 0.9690761713059198
$ws.Range("G11").Value = 1.060223973180004
$ws.Range("H11").Value = 0.8554390596845817
$ws.Range("I11").Value = 1.039363489602529
$ws.Range("J11").Value = 0.9738239194210393
$ws.Range("K11").Value = 1.02313933809295
$ws.Range("L11").Value = 0.9107500069952641
$ws.Range("M11").Value = 1.060223973180004
$ws.Range("N11").Value = 0.9422631034583291
$ws.Range("O11").Value = 0.9784565878506455
$ws.Range("P11").Value = 0.9826128881892956

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 1.058037579436627
$ws.Range("D12").Value = 0.8570081536378086
$ws.Range("E12").Value = 1.029149849972901
$ws.Range("F12").Value = 0.9695497690206656
$ws.Range("G12").Value = 1.058037579436627
$ws.Range("H12").Value = 0.8570081536378086
$ws.Range("I12").Value = 1.038759582526452
$ws.Range("J12").Value = 0.9745227022647185
$ws.Range("K12").Value = 1.022406277314573
$ws.Range("L12").Value = 0.911997874409501
$ws.Range("M12").Value = 1.058037579436627
$ws.Range("N12").Value = 0.9430790018053548
$ws.Range("O12").Value = 0.9784363380170007
$ws.Range("P12").Value = 0.9826789735729059

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 1.059851375813614
$ws.Range("D13").Value = 0.8558534750528169
$ws.Range("E13").Value = 1.02895000138764
$ws.Range("F13").Value = 0.9692219342567756
$ws.Range("G13").Value = 1.059851375813614
$ws.Range("H13").Value = 0.8558534750528169
$ws.Range("I13").Value = 1.039169011503091
$ws.Range("J13").Value = 0.9739366715468677
$ws.Range("K13").Value = 1.022993049123258
$ws.Range("L13").Value = 0.9110816282283328
$ws.Range("M13").Value = 1.059851375813614
$ws.Range("N13").Value = 0.9424017382202285
$ws.Range("O13").Value = 0.9784691966277118
$ws.Range("P13").Value = 0.9826321433640497

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.8969480000000002
$ws.Range("D14").Value = 0.7726800000000016
$ws.Range("E14").Value = 1.211807999999999
$ws.Range("F14").Value = 0.9403520000000002
$ws.Range("G14").Value = 0.8969480000000002
$ws.Range("H14").Value = 0.7726800000000016
$ws.Range("I14").Value = 1.097607999999999
$ws.Range("J14").Value = 1.038947999999999
$ws.Range("K14").Value = 0.9683600000000006
$ws.Range("L14").Value = 0.8762999999999989
$ws.Range("M14").Value = 0.8969480000000002
$ws.Range("N14").Value = 0.9922440000000003
$ws.Range("O14").Value = 0.9554470000000003
$ws.Range("P14").Value = 0.9753755

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.96
$ws.Range("D15").Value = 0.53
$ws.Range("E15").Value = 1.37
$ws.Range("F15").Value = 0.87
$ws.Range("G15").Value = 0.96
$ws.Range("H15").Value = 0.53
$ws.Range("I15").Value = 1.21
$ws.Range("J15").Value = 1.025962499999997
$ws.Range("K15").Value = 0.99
$ws.Range("L15").Value = 0.72
$ws.Range("M15").Value = 0.96
$ws.Range("N15").Value = 0.9500000000000001
$ws.Range("O15").Value = 0.9325000000000001
$ws.Range("P15").Value = 0.9594953124999998

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.9776814460927974
$ws.Range("D16").Value = 0.7236039903232024
$ws.Range("E16").Value = 1.2101315629056
$ws.Range("F16").Value = 0.9229775804416026
$ws.Range("G16").Value = 0.9776814460927974
$ws.Range("H16").Value = 0.7236039903232024
$ws.Range("I16").Value = 1.117745667891196
$ws.Range("J16").Value = 1.010792811315193
$ws.Range("K16").Value = 0.9926541283328028
$ws.Range("L16").Value = 0.8351458914304004
$ws.Range("M16").Value = 0.9776790691839974
$ws.Range("N16").Value = 0.9668677766144014
$ws.Range("O16").Value = 0.9585986449408007
$ws.Range("P16").Value = 0.9738416348415995

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9954657859163153
$ws.Range("D17").Value = 0.9946426984781939
$ws.Range("E17").Value = 0.9926040042299493
$ws.Range("F17").Value = 0.9933143444385427
$ws.Range("G17").Value = 0.9954657859163153
$ws.Range("H17").Value = 0.9946426984781939
$ws.Range("I17").Value = 0.9933958695919379
$ws.Range("J17").Value = 0.9936341623097674
$ws.Range("K17").Value = 0.9940882974342137
$ws.Range("L17").Value = 0.9930312882925031
$ws.Range("M17").Value = 0.9954573140957049
$ws.Range("N17").Value = 0.9936233513540715
$ws.Range("O17").Value = 0.9940067082657503
$ws.Range("P17").Value = 0.9937720563364278

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9695898464352465
$ws.Range("D18").Value = 0.9960400781831708
$ws.Range("E18").Value = 1.004756862539783
$ws.Range("F18").Value = 0.9948983831330787
$ws.Range("G18").Value = 0.9695898464352465
$ws.Range("H18").Value = 0.9960400781831708
$ws.Range("I18").Value = 0.9940099476001841
$ws.Range("J18").Value = 1.004051893565846
$ws.Range("K18").Value = 0.9873311808283504
$ws.Range("L18").Value = 0.9969009270266295
$ws.Range("M18").Value = 0.9695898464352465
$ws.Range("N18").Value = 1.000398470361477
$ws.Range("O18").Value = 0.9913212925728198
$ws.Range("P18").Value = 0.9934473899140361

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9728325754911196
$ws.Range("D19").Value = 1.034135537007103
$ws.Range("E19").Value = 0.984185010272204
$ws.Range("F19").Value = 1.00119227830787
$ws.Range("G19").Value = 0.9728325754911196
$ws.Range("H19").Value = 1.034135537007103
$ws.Range("I19").Value = 0.9781234052084472
$ws.Range("J19").Value = 0.9985072133579688
$ws.Range("K19").Value = 0.9846227626898721
$ws.Range("L19").Value = 1.019178046933106
$ws.Range("M19").Value = 0.9728155426728395
$ws.Range("N19").Value = 1.009160273639653
$ws.Range("O19").Value = 0.998086350269574
$ws.Range("P19").Value = 0.9965971036584613

# Apply the same formatting (style) used by the other row-index cells in column A
# to the three newly added rows (17-19).
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
